# This edit reshuffles the data rows (2-37) of the sheet: each target row's
# D and K:T values are taken from a different source row (a pure
# permutation of the existing rows -- same data, new order). Column D is a
# date (serial number) that carries a date-formatted style already applied
# to the cell, so only .Value2 needs to be updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row number -> source row number (the row whose data, prior to
# this edit, should end up in the target row).
$map = @{
    2  = 8
    3  = 9
    4  = 25
    5  = 23
    6  = 30
    7  = 17
    8  = 10
    9  = 11
    10 = 20
    11 = 24
    12 = 6
    13 = 7
    14 = 36
    15 = 37
    16 = 2
    17 = 3
    18 = 4
    19 = 5
    20 = 21
    21 = 22
    22 = 18
    23 = 26
    24 = 27
    25 = 28
    26 = 29
    27 = 16
    28 = 31
    29 = 13
    30 = 14
    31 = 15
    32 = 34
    33 = 35
    34 = 32
    35 = 33
    36 = 12
    37 = 19
}

# Columns whose values move together with each row.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the "before" values for every row/column we might read from,
# since several writes will overwrite rows that are later used as sources.
$snapshot = @{}
for ($r = 2; $r -le 37; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the snapshotted source-row data into each target row.
foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $srcVals[$col]
    }
}
